$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cell values for columns A and B, rows 1-32 (per diff)
$ws.Range("A1").Value = -0.092826440534025778
$ws.Range("B1").Value = 0.092701678350678662
$ws.Range("A2").Value = 0.0012287138923614549
$ws.Range("B2").Value = -0.0017718058254860836
$ws.Range("A3").Value = 0.10554578879861864
$ws.Range("B3").Value = -0.10596545026620063
$ws.Range("A4").Value = -0.1780311985376315
$ws.Range("B4").Value = 0.1771367691936554
$ws.Range("A5").Value = -0.17113676938005806
$ws.Range("B5").Value = 0.16934434904417195
$ws.Range("A6").Value = -0.095578244970149662
$ws.Range("B6").Value = 0.09545330253978257
$ws.Range("A7").Value = -0.075453302766618435
$ws.Range("B7").Value = 0.07516377864235757
$ws.Range("A8").Value = -0.055163778871374802
$ws.Range("B8").Value = 0.054945556012992114
$ws.Range("A9").Value = -0.048945556211484664
$ws.Range("B9").Value = 0.048770910545306023
$ws.Range("A10").Value = -0.042770910746043
$ws.Range("B10").Value = 0.042746989207032016
$ws.Range("A11").Value = -0.038246989404260745
$ws.Range("B11").Value = 0.038210323829382276
$ws.Range("A12").Value = -0.032210324030724102
$ws.Range("B12").Value = 0.032107562486773933
$ws.Range("A13").Value = -0.026107562690269148
$ws.Range("B13").Value = 0.026082806341202414
$ws.Range("A14").Value = -0.014082806559489569
$ws.Range("B14").Value = 0.014072107697031377
$ws.Range("A15").Value = -0.0080721079014569597
$ws.Range("B15").Value = 0.0080642579935350511
$ws.Range("A16").Value = -0.002064258198202662
$ws.Range("B16").Value = 0.002051819690349177
$ws.Range("A17").Value = -0.0090046006164010706
$ws.Range("B17").Value = 0.0089999997879406379
$ws.Range("A18").Value = -0.083930392437249424
$ws.Range("B18").Value = 0.083854725218579063
$ws.Range("A19").Value = -0.027096972347577353
$ws.Range("B19").Value = 0.027014065757839045
$ws.Range("A20").Value = -0.018014065951822644
$ws.Range("B20").Value = 0.018004308536403002
$ws.Range("A21").Value = -0.0090043087306295178
$ws.Range("B21").Value = 0.0089999998056011776
$ws.Range("A22").Value = -0.094773983165174158
$ws.Range("B22").Value = 0.094465614726878755
$ws.Range("A23").Value = -0.085465614920549271
$ws.Range("B23").Value = 0.084966252751042859
$ws.Range("A24").Value = -0.042966253025905843
$ws.Range("B24").Value = 0.042841264841605664
$ws.Range("A25").Value = -0.067766104266929972
$ws.Range("B25").Value = 0.06769105038926071
$ws.Range("A26").Value = -0.061691050582663109
$ws.Range("B26").Value = 0.061602250838880934
$ws.Range("A27").Value = -0.055602251032808248
$ws.Range("B27").Value = 0.055326337939896497
$ws.Range("A28").Value = -0.049326338136476799
$ws.Range("B28").Value = 0.049156302286211151
$ws.Range("A29").Value = -0.063009230961254659
$ws.Range("B29").Value = 0.062787380935940718
$ws.Range("A30").Value = -0.042787381168606142
$ws.Range("B30").Value = 0.042663244261443189
$ws.Range("A31").Value = -0.027663244484134708
$ws.Range("B31").Value = 0.027642477481096606
$ws.Range("A32").Value = -0.0066424777183007322
$ws.Range("B32").Value = 0.0066262680789419548

# Widen column B to match column A (15.42578125); the host quantizes
# ColumnWidth to 1/6-character steps, so 14.67 is the closest input that
# lands on the nearest achievable stored width to the target.
$ws.Columns.Item(2).ColumnWidth = 14.67
